$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 729.2368
$ws.Range("J17").Value = 747.5278
$ws.Range("L17").Value = 2242.5834
$ws.Range("N17").Value = -2578.5834

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 480971.2
$ws.Range("I107").Value = 635272.9
$ws.Range("J107").Value = 921.55554
$ws.Range("K107").Value = 635272.9
$ws.Range("L107").Value = 921.55554
$ws.Range("M107").Value = -633352.9
$ws.Range("N107").Value = -4761.55554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 28514

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2607.6
$ws.Range("I45").Value = 1410
$ws.Range("K45").Value = 1410
$ws.Range("M45").Value = -1033

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2076.5557
$ws.Range("I61").Value = 1675.4736
$ws.Range("J61").Value = 3029.125
$ws.Range("K61").Value = 1675.4736
$ws.Range("L61").Value = 3029.125
$ws.Range("M61").Value = -1463.4736
$ws.Range("N61").Value = -3453.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2864.5334
$ws.Range("I122").Value = 3211.25
$ws.Range("J122").Value = 2468.2856
$ws.Range("K122").Value = 9633.75
$ws.Range("L122").Value = 7404.8568
$ws.Range("M122").Value = -7183.75
$ws.Range("N122").Value = -12304.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4079.2222
$ws.Range("I132").Value = 3449
$ws.Range("J132").Value = 4583.4
$ws.Range("K132").Value = 10347
$ws.Range("L132").Value = 13750.2
$ws.Range("M132").Value = -7817
$ws.Range("N132").Value = -18810.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2076.5557
$ws.Range("I136").Value = 1675.4736
$ws.Range("J136").Value = 3029.125
$ws.Range("K136").Value = 5026.4208
$ws.Range("L136").Value = 9087.375
$ws.Range("M136").Value = -2476.4208
$ws.Range("N136").Value = -14187.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 7375.25
$ws.Range("I5").Value = 8214.571
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 8214.571
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = -8101.571
$ws.Range("N5").Value = -1726

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value = 26111
$ws.Range("J116").Value = 26111
$ws.Range("L116").Value = 26111
$ws.Range("N116").Value = -35289

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2203.1702
$ws.Range("I134").Value = 1286.5
$ws.Range("K134").Value = 3859.5
$ws.Range("M134").Value = -1324.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 29462.625
$ws.Range("I62").Value = 42379
$ws.Range("J62").Value = 7935.3335
$ws.Range("K62").Value = 42379
$ws.Range("L62").Value = 7935.3335
$ws.Range("M62").Value = -41755
$ws.Range("N62").Value = -9183.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 29462.625
$ws.Range("I65").Value = 42379
$ws.Range("J65").Value = 7935.3335
$ws.Range("K65").Value = 211895
$ws.Range("L65").Value = 39676.6675
$ws.Range("M65").Value = -208775
$ws.Range("N65").Value = -45916.6675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2509.9697
$ws.Range("I132").Value = 1381.8182
$ws.Range("K132").Value = 4145.4546
$ws.Range("M132").Value = -1615.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1897.1666
$ws.Range("I134").Value = 727.3333
$ws.Range("J134").Value = 4626.778
$ws.Range("K134").Value = 2181.9999
$ws.Range("L134").Value = 13880.334
$ws.Range("M134").Value = 353.0001000000002
$ws.Range("N134").Value = -18950.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 27907.455
$ws.Range("I4").Value = 245.5
$ws.Range("J4").Value = 43714.285
$ws.Range("K4").Value = 736.5
$ws.Range("L4").Value = 131142.855
$ws.Range("M4").Value = -624.5
$ws.Range("N4").Value = -131366.855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2478.8735
$ws.Range("I68").Value = 3279.0208
$ws.Range("K68").Value = 9837.062399999999
$ws.Range("M68").Value = -9026.062399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2478.8735
$ws.Range("I71").Value = 3279.0208
$ws.Range("K71").Value = 29511.1872
$ws.Range("M71").Value = -25455.1872

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 834823.0600000001
$ws.Range("I107").Value = 322.25
$ws.Range("J107").Value = 1252073.5
$ws.Range("K107").Value = 966.75
$ws.Range("L107").Value = 3756220.5
$ws.Range("M107").Value = 953.25
$ws.Range("N107").Value = -3760060.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 3035.3572
$ws.Range("I108").Value = 2541.25
$ws.Range("K108").Value = 7623.75
$ws.Range("M108").Value = -4743.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2410.9092
$ws.Range("I109").Value = 1250
$ws.Range("J109").Value = 2668.889
$ws.Range("K109").Value = 3750
$ws.Range("L109").Value = 8006.667
$ws.Range("M109").Value = -2710
$ws.Range("N109").Value = -10086.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 399.6111
$ws.Range("J113").Value = 416.20834
$ws.Range("L113").Value = 1248.62502
$ws.Range("N113").Value = -5588.625019999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2746.6196
$ws.Range("J131").Value = 3123.5737
$ws.Range("L131").Value = 9370.721099999999
$ws.Range("N131").Value = -19450.7211

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 50002.668
$ws.Range("I4").Value = 10000
$ws.Range("K4").Value = 10000
$ws.Range("M4").Value = -9888

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2224642.2
$ws.Range("J122").Value = 3200
$ws.Range("L122").Value = 9600
$ws.Range("N122").Value = -14500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2698.175
$ws.Range("I132").Value = 2218.7
$ws.Range("K132").Value = 6656.099999999999
$ws.Range("M132").Value = -4126.099999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3228.9167
$ws.Range("I7").Value = 2415.6667
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 2415.6667
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -2303.6667
$ws.Range("N7").Value = -3724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 27000.5
$ws.Range("I39").Value = 27000.5
$ws.Range("K39").Value = 27000.5
$ws.Range("M39").Value = -26540.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3387.1738
$ws.Range("I122").Value = 2475
$ws.Range("J122").Value = 3579.2104
$ws.Range("K122").Value = 7425
$ws.Range("L122").Value = 10737.6312
$ws.Range("M122").Value = -4975
$ws.Range("N122").Value = -15637.6312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3228.9167
$ws.Range("I126").Value = 2415.6667
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 7247.000100000001
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -4777.000100000001
$ws.Range("N126").Value = -15440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3348.6943
$ws.Range("I132").Value = 2507.64
$ws.Range("K132").Value = 7522.92
$ws.Range("M132").Value = -4992.92

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 7386.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 45245.086
$ws.Range("I122").Value = 112392.445
$ws.Range("J122").Value = 2078.9285
$ws.Range("K122").Value = 337177.335
$ws.Range("L122").Value = 6236.7855
$ws.Range("M122").Value = -334727.335
$ws.Range("N122").Value = -11136.7855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11629914
$ws.Range("I132").Value = 15626507
$ws.Range("J132").Value = 3463.6365
$ws.Range("K132").Value = 46879521
$ws.Range("L132").Value = 10390.9095
$ws.Range("M132").Value = -46876991
$ws.Range("N132").Value = -15450.9095
